$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (CIVN2020 conference entry) now links the congress name, and the row
# grows taller to fit the wrapped link text.
$ws.Cells.Item(2, 3).Value = "\href{https://www.youtube.com/playlist?list=PLI4QwBEXHFJKVhSBqaof6gJdKCi3CN2UM}{CIVN2020 - Congreso Internacional de Neurociencias: Cerebro y Comportamiento en Tiempos de COVID-19}"
$ws.Rows.Item(2).RowHeight = 75

# Row 4 (ISEP congress entry) now links the congress name, and picks up an
# explicit row height.
$ws.Cells.Item(4, 3).Value = "\href{https://www.isep.es/congreso2020/}{1er Congreso Internacional Virtual ISEP}"
$ws.Rows.Item(4).RowHeight = 30

# Selection moves from B10 to C10.
$ws.Range("C10").Select()
